{"js": "// The template's field codes (e.g. \" m: 2.myTemplate() \") used to be stored\n// as real Word fields (fldChar begin / instrText / fldChar end). The parser\n// was updated (TokenIteratorFieldRewriterSplit) to read the tag straight out\n// of literal run text instead, so every `m:` field in the document is\n// rewritten here as plain text \"{<trimmed field code>}\" in place of the\n// field, keeping the paragraph's run formatting (language).\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nfor (const paragraph of paragraphs.items) {\n  // Capture the paragraph's own range/language before any edits so the\n  // replacement text run keeps the same formatting as the field runs had.\n  const paragraphRange = paragraph.getRange();\n  paragraphRange.load(\"languageId\");\n\n  const paragraphFields = paragraph.fields;\n  paragraphFields.load(\"items\");\n  await context.sync();\n\n  if (paragraphFields.items.length === 0) {\n    continue;\n  }\n\n  for (const field of paragraphFields.items) {\n    field.load(\"code\");\n  }\n  await context.sync();\n\n  for (const field of paragraphFields.items) {\n    // \" m: 2.myTemplate() \" -> \"m: 2.myTemplate()\"\n    const code = field.code.trim();\n    const languageId = paragraphRange.languageId;\n\n    // Result is the (empty) range right where the field lived; grab it\n    // before deleting the field so we can drop the literal text there.\n    const resultRange = field.result;\n    field.delete();\n    await context.sync();\n\n    const insertedRange = resultRange.insertText(\"{\" + code + \"}\", Word.InsertLocation.replace);\n    insertedRange.languageId = languageId;\n    await context.sync();\n  }\n}\n", "ps1": "# The template's field codes (e.g. \" m: 2.myTemplate() \") were previously\n# stored as real Word fields (fldChar begin / instrText / fldChar end).\n# The parser was updated (TokenIteratorFieldRewriterSplit) to instead read\n# the tag straight out of literal run text, so every `m:` field in the\n# document is rewritten here as plain text \"{<trimmed field code>}\" in\n# place of the field, keeping the paragraph's run formatting (language).\n\n$d = $word.ActiveDocument\n\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $para = $d.Paragraphs.Item($i)\n    $paraFields = $para.Range.Fields\n\n    while ($paraFields.Count -gt 0) {\n        $field = $paraFields.Item(1)\n\n        # Field code text, e.g. \" m: 2.myTemplate() \" -> \"m: 2.myTemplate()\"\n        $code = $field.Code.Text.Trim()\n\n        # Grab the paragraph's language before we touch anything, so the\n        # replacement run keeps the same formatting as the deleted field runs.\n        $paraRange = $para.Range\n        $langId = $paraRange.LanguageID\n\n        # Removes the begin/instrText/end runs that made up the field.\n        $field.Delete()\n\n        # Re-fetch the (now field-less) paragraph range and drop the plain\n        # \"{...}\" text in where the field used to be.\n        $insertRange = $d.Paragraphs.Item($i).Range\n        $insertRange.InsertBefore(\"{\" + $code + \"}\")\n        $insertRange.LanguageID = $langId\n\n        $paraFields = $para.Range.Fields\n    }\n}\n"}
